$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest 59 minutes of timestamp (col A) / ETH price (col B) data from the
# scalping-bot's live feed, replacing the previous rolling window.
$aValues = @(0.17847222222222223,0.17916666666666667,0.17986111111111111,0.18055555555555555,0.18124999999999999,0.18194444444444444,0.18263888888888891,0.18333333333333335,0.18402777777777779,0.18472222222222223,0.18541666666666667,0.18611111111111112,0.18680555555555556,0.1875,0.18819444444444444,0.18888888888888888,0.18958333333333333,0.19027777777777777,0.19097222222222221,0.19166666666666665,0.19236111111111112,0.19305555555555554,0.19375000000000001,0.19444444444444445,0.19513888888888889,0.19583333333333333,0.19652777777777777,0.19722222222222222,0.19791666666666666,0.1986111111111111,0.19930555555555554,0.19999999999999998,0.20069444444444443,0.20138888888888887,0.20208333333333331,0.20277777777777781,0.20347222222222219,0.20416666666666669,0.20486111111111113,0.20555555555555557,0.20625000000000002,0.20694444444444446,0.2076388888888889,0.20833333333333334,0.20902777777777778,0.20972222222222223,0.21041666666666667,0.21111111111111111,0.21180555555555555,0.21249999999999999,0.21319444444444444,0.21388888888888891,0.21458333333333335,0.21527777777777779,0.21597222222222223,0.21666666666666667,0.21736111111111112,0.21805555555555556,0.21875)

$bValues = @(1584.1458334184999,1582.8766285316999,1582.6705184717,1584.0122182032001,1585.7084572923,1585.1022675858001,1585.0897337214001,1583.9639754264999,1583.6197433350001,1584.7136449728,1584.7200664612999,1584.0820353697,1584.2590199069,1584.1921882951999,1583.9392918282999,1584.6874661001,1583.9199292901001,1583.3895764177,1583.0611068170999,1582.3049278596,1583.2786001125,1583.2411580884,1583.2906463694001,1584.1736097116,1584.290889201,1584.5652183489999,1584.2681410651001,1584.7298111053001,1584.1510144275001,1584.1930937606,1583.8388058343,1583.5974854336,1583.2473179131,1582.8619468493,1583.1919116289,1583.2028641366001,1583.0516249036,1583.2223956093001,1583.5474802341,1583.3768332297,1583.0432250895001,1583.5837492488999,1583.4422379969999,1583.150541216,1583.8502882800001,1583.7752729424001,1583.7585477335001,1584.198046774,1584.0052795588001,1583.8094178346,1583.8346487558999,1584.086560511,1584.3370469035001,1584.0889417763001,1584.0250986998999,1584.2879790161001,1584.2418544309,1584.1907570047999,1584.0853161048001)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

$ws.Range("B1:B59").Select() | Out-Null
